# Updates cryptos list values (prices/volumes) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.805.60'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '3.360.67'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  +0.03%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.63'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -0.16%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.91'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('E10').Value = '  -2.48%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.380'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  -3.75%  '
$ws.Range('D12').Value = '3.934.35'
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('E13').Value = '  +1.67%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.54'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').Value = '3.357.18'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('E16').Value = '  -2.66%  '
$ws.Range('D17').Value = '60.897.53'
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.05'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  -3.37%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.50'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  -4.05%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.81'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -2.42%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.81'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -1.56%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.21'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +1.89%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.546'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -2.50%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('E25').Value = '  -6.00%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.189'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +6.19%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  -4.13%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.80'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('E31').Value = '  -0.06%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.34'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -6.05%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.83'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  -3.66%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.83'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  -1.80%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '165.94'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -0.81%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.90'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('D37').Value = '3.396.89'
$ws.Range('E37').Value = '  -0.61%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.44'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  -4.06%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0760'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '25.19'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  -10.55%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.770'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('E42').Value = '  -2.64%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.62'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -3.86%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.10'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('D45').Value = '2.445.49'
$ws.Range('E45').Value = '  -4.19%  '
$ws.Range('E46').Value = '  +0.03%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.57'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -4.08%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.07'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  -5.90%  '
$ws.Range('E49').Value = '  -5.03%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.02'
$ws.Range('D50').Style = $origStyle
$ws.Range('E51').Value = '  -3.54%  '
